$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and Litecoin/PEPE row swap)
# from the latest GitHub Actions scrape. Values are prefixed with a
# leading apostrophe so Excel stores them as literal text (matching
# the original inline-string cell contents) instead of auto-converting
# number-like strings (e.g. "646.74", "1.00") into numeric values.

# Row 2
$ws.Range("D2").Value = "'68.915.46"
$ws.Range("E2").Value = "'  -1.07%  "
# Row 3
$ws.Range("D3").Value = "'3.756.78"
$ws.Range("E3").Value = "'  -2.01%  "
# Row 4
$ws.Range("E4").Value = "'  -0.14%  "
# Row 5
$ws.Range("D5").Value = "'646.74"
$ws.Range("E5").Value = "'  +2.27%  "
# Row 6
$ws.Range("D6").Value = "'165.52"
$ws.Range("E6").Value = "'  -0.47%  "
# Row 7
$ws.Range("D7").Value = "'3.753.28"
$ws.Range("E7").Value = "'  -2.01%  "
# Row 8
$ws.Range("E8").Value = "'  +0.10%  "
# Row 9
$ws.Range("E9").Value = "'  +0.31%  "
# Row 10
$ws.Range("E10").Value = "'  -2.28%  "
# Row 11
$ws.Range("E11").Value = "'  +0.17%  "
# Row 12
$ws.Range("D12").Value = "'6.92"
$ws.Range("E12").Value = "'  +3.63%  "
# Row 13
$ws.Range("E13").Value = "'  -5.13%  "
# Row 14
$ws.Range("D14").Value = "'34.84"
$ws.Range("E14").Value = "'  -3.55%  "
# Row 15
$ws.Range("D15").Value = "'4.391.04"
$ws.Range("E15").Value = "'  -1.83%  "
# Row 16
$ws.Range("D16").Value = "'3.756.68"
$ws.Range("E16").Value = "'  -3.16%  "
# Row 17
$ws.Range("D17").Value = "'68.923.67"
$ws.Range("E17").Value = "'  -0.96%  "
# Row 18
$ws.Range("D18").Value = "'17.72"
$ws.Range("E18").Value = "'  -2.76%  "
# Row 19
$ws.Range("E19").Value = "'  +0.08%  "
# Row 20
$ws.Range("D20").Value = "'6.99"
$ws.Range("E20").Value = "'  -2.44%  "
# Row 21
$ws.Range("D21").Value = "'466.36"
$ws.Range("E21").Value = "'  -0.55%  "
# Row 22
$ws.Range("D22").Value = "'9.53"
$ws.Range("E22").Value = "'  -1.95%  "
# Row 23
$ws.Range("E23").Value = "'  -0.68%  "
# Row 24
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'81.55"
$ws.Range("E24").Value = "'  -2.77%  "
# Row 25
$ws.Range("B25").Value = "'PEPE"
$ws.Range("C25").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000143"
$ws.Range("E25").Value = "'  -6.40%  "
# Row 26
$ws.Range("D26").Value = "'12.21"
$ws.Range("E26").Value = "'  +1.46%  "
# Row 27
$ws.Range("D27").Value = "'10.13"
$ws.Range("E27").Value = "'  +0.19%  "
# Row 28
$ws.Range("D28").Value = "'2.09"
$ws.Range("E28").Value = "'  -4.57%  "
# Row 29
$ws.Range("E29").Value = "'  -0.09%  "
# Row 30
$ws.Range("D30").Value = "'3.903.49"
$ws.Range("E30").Value = "'  -1.91%  "
# Row 31
$ws.Range("E31").Value = "'  -0.35%  "
# Row 32
$ws.Range("E32").Value = "'  +1.18%  "
# Row 33
$ws.Range("D33").Value = "'7.12"
$ws.Range("E33").Value = "'  -2.64%  "
# Row 34
$ws.Range("D34").Value = "'28.52"
$ws.Range("E34").Value = "'  -2.81%  "
# Row 35
$ws.Range("D35").Value = "'0.173"
$ws.Range("E35").Value = "'  +15.31%  "
# Row 36
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  -0.34%  "
# Row 37
$ws.Range("D37").Value = "'3.710.31"
$ws.Range("E37").Value = "'  -1.75%  "
# Row 38
$ws.Range("D38").Value = "'8.79"
$ws.Range("E38").Value = "'  -3.19%  "
# Row 39
$ws.Range("E39").Value = "'  -3.35%  "
# Row 40
$ws.Range("E40").Value = "'  -2.93%  "
# Row 41
$ws.Range("E41").Value = "'  -7.18%  "
# Row 42
$ws.Range("E42").Value = "'  +0.18%  "
# Row 43
$ws.Range("D43").Value = "'0.954"
$ws.Range("E43").Value = "'  -2.96%  "
# Row 44
$ws.Range("E44").Value = "'  -0.02%  "
# Row 45
$ws.Range("D45").Value = "'45.04"
$ws.Range("E45").Value = "'  +2.76%  "
# Row 46
$ws.Range("D46").Value = "'1.97"
$ws.Range("E46").Value = "'  +1.93%  "
# Row 47
$ws.Range("D47").Value = "'154.84"
$ws.Range("E47").Value = "'  -0.83%  "
# Row 48
$ws.Range("D48").Value = "'47.15"
$ws.Range("E48").Value = "'  +0.25%  "
# Row 49
$ws.Range("D49").Value = "'0.294"
$ws.Range("E49").Value = "'  -2.91%  "
# Row 50
$ws.Range("E50").Value = "'  -1.18%  "
# Row 51
$ws.Range("D51").Value = "'8.34"
$ws.Range("E51").Value = "'  -1.60%  "
